# Generate Report for Handback
# This script updates the localization-status workbook to reflect a completed
# handback: status text changes, new handback timestamps, newly populated
# "Latest Target File" / "Latest Handback File" columns (with hyperlinks) on
# the per-locale sheets, and widened columns to accommodate the longer values.

$wb = $excel.ActiveWorkbook

$colWidth40 = 39.166666666666664   # renders as stored width 40
$colWidthWide = 29.166666666666664 # renders as stored width ~30 (target 29.9777050018311)

# ---------------------------------------------------------------------------
# Overview sheet: status text ("Ready for handoff" -> "Handed back: in sync
# with en-US") is a shared string referenced from E2:F3 as well, so updating
# those cells (and the matching cells on the locale sheets) updates it
# everywhere it is used. Also widen columns E/F to fit the new, longer text.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Range("E1").ColumnWidth = $colWidthWide
$wsOverview.Range("F1").ColumnWidth = $colWidthWide

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column (same shared string as Overview E/F)
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

# Latest Target File (I) / Latest Handback File (J) are now populated
$wsZh.Range("I2").Value = "c3b7d116-9e01-4b24-b951-e98c9c14e9aa.md"
$wsZh.Range("J2").Value = "c3b7d116-9e01-4b24-b951-e98c9c14e9aa.25237311146c3ae68ca7f93139b46d06625f134a.zh-cn.xlf"
$wsZh.Range("I3").Value = "ec4dc51d-aefd-4ed3-ba40-28d3c535c121.md"
$wsZh.Range("J3").Value = "ec4dc51d-aefd-4ed3-ba40-28d3c535c121.87cc1f4c8a38acf00a3afa9da2d9edc8abce26ac.zh-cn.xlf"

# Latest Handback DateTime (K) now has a real timestamp
$wsZh.Range("K2").Value = "2016-09-06 16:59:33"
$wsZh.Range("K3").Value = "2016-09-06 16:59:33"

# Hyperlinks for the newly populated Latest Target File cells
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84b83d213f48e0896ca5aa5813d31891b4b7a000/e2e/c3b7d116-9e01-4b24-b951-e98c9c14e9aa.md", "", "", "c3b7d116-9e01-4b24-b951-e98c9c14e9aa.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84b83d213f48e0896ca5aa5813d31891b4b7a000/e2e/ec4dc51d-aefd-4ed3-ba40-28d3c535c121.md", "", "", "ec4dc51d-aefd-4ed3-ba40-28d3c535c121.md")

$wsZh.Range("C1").ColumnWidth = $colWidthWide
$wsZh.Range("I1").ColumnWidth = $colWidth40
$wsZh.Range("J1").ColumnWidth = $colWidth40

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column (same shared string as Overview E/F)
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

# Latest Target File (I) / Latest Handback File (J) are now populated
$wsDe.Range("I2").Value = "c3b7d116-9e01-4b24-b951-e98c9c14e9aa.md"
$wsDe.Range("J2").Value = "c3b7d116-9e01-4b24-b951-e98c9c14e9aa.25237311146c3ae68ca7f93139b46d06625f134a.de-de.xlf"
$wsDe.Range("I3").Value = "ec4dc51d-aefd-4ed3-ba40-28d3c535c121.md"
$wsDe.Range("J3").Value = "ec4dc51d-aefd-4ed3-ba40-28d3c535c121.87cc1f4c8a38acf00a3afa9da2d9edc8abce26ac.de-de.xlf"

# Latest Handback DateTime (K) now has a real timestamp
$wsDe.Range("K2").Value = "2016-09-06 16:59:52"
$wsDe.Range("K3").Value = "2016-09-06 16:59:52"

# Hyperlinks for the newly populated Latest Target File cells
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84b83d213f48e0896ca5aa5813d31891b4b7a000/e2e/c3b7d116-9e01-4b24-b951-e98c9c14e9aa.md", "", "", "c3b7d116-9e01-4b24-b951-e98c9c14e9aa.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84b83d213f48e0896ca5aa5813d31891b4b7a000/e2e/ec4dc51d-aefd-4ed3-ba40-28d3c535c121.md", "", "", "ec4dc51d-aefd-4ed3-ba40-28d3c535c121.md")

$wsDe.Range("C1").ColumnWidth = $colWidthWide
$wsDe.Range("I1").ColumnWidth = $colWidth40
$wsDe.Range("J1").ColumnWidth = $colWidth40
